$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A66").Value = "MF49C0012"
$ws.Range("A67").Value = "MF49F0062"
$ws.Range("A68").Value = "MF4970266"
$ws.Range("A69").Value = "MF49G0263"
